$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "ba"
$ws.Range("J3").Value = "Appreciation"
$ws.Range("I10").Value = "sv"
$ws.Range("J10").Value = "Statement-opinion"
$ws.Range("I12").Value = "%"
$ws.Range("J12").Value = "Uninterpretable"
$ws.Range("I21").Value = "sv"
$ws.Range("J21").Value = "Statement-opinion"
$ws.Range("I27").Value = "aa"
$ws.Range("J27").Value = "Agree/Accept"
$ws.Range("I45").Value = "sv"
$ws.Range("J45").Value = "Statement-opinion"
$ws.Range("I52").Value = "sd"
$ws.Range("J52").Value = "Statement-non-opinion"
$ws.Range("I60").Value = "sv"
$ws.Range("J60").Value = "Statement-opinion"
$ws.Range("I68").Value = "aa"
$ws.Range("J68").Value = "Agree/Accept"
$ws.Range("I72").Value = "sv"
$ws.Range("J72").Value = "Statement-opinion"
$ws.Range("I108").Value = "sd"
$ws.Range("J108").Value = "Statement-non-opinion"
$ws.Range("I127").Value = "sv"
$ws.Range("J127").Value = "Statement-opinion"
$ws.Range("I134").Value = "sd"
$ws.Range("J134").Value = "Statement-non-opinion"
$ws.Range("I158").Value = "sv"
$ws.Range("J158").Value = "Statement-opinion"
$ws.Range("I163").Value = "sv"
$ws.Range("J163").Value = "Statement-opinion"
$ws.Range("I175").Value = "b"
$ws.Range("J175").Value = "Acknowledge (Backchannel)"
$ws.Range("I180").Value = "sd"
$ws.Range("J180").Value = "Statement-non-opinion"
$ws.Range("I203").Value = "sv"
$ws.Range("J203").Value = "Statement-opinion"
$ws.Range("I212").Value = "sd"
$ws.Range("J212").Value = "Statement-non-opinion"
$ws.Range("I214").Value = "sd"
$ws.Range("J214").Value = "Statement-non-opinion"
$ws.Range("I221").Value = "ba"
$ws.Range("J221").Value = "Appreciation"
$ws.Range("I224").Value = "sv"
$ws.Range("J224").Value = "Statement-opinion"
$ws.Range("I230").Value = "b"
$ws.Range("J230").Value = "Acknowledge (Backchannel)"
$ws.Range("I246").Value = "aa"
$ws.Range("J246").Value = "Agree/Accept"
$ws.Range("I247").Value = "aa"
$ws.Range("J247").Value = "Agree/Accept"
$ws.Range("I248").Value = "%"
$ws.Range("J248").Value = "Uninterpretable"
$ws.Range("I258").Value = "aa"
$ws.Range("J258").Value = "Agree/Accept"
$ws.Range("I277").Value = "aa"
$ws.Range("J277").Value = "Agree/Accept"
$ws.Range("I285").Value = "sd"
$ws.Range("J285").Value = "Statement-non-opinion"
$ws.Range("I292").Value = "sd"
$ws.Range("J292").Value = "Statement-non-opinion"
$ws.Range("I293").Value = "sv"
$ws.Range("J293").Value = "Statement-opinion"
$ws.Range("I296").Value = "sd"
$ws.Range("J296").Value = "Statement-non-opinion"
$ws.Range("I310").Value = "sd"
$ws.Range("J310").Value = "Statement-non-opinion"
$ws.Range("I327").Value = "b"
$ws.Range("J327").Value = "Acknowledge (Backchannel)"
$ws.Range("I344").Value = "qy"
$ws.Range("J344").Value = "Yes-No-Question"
$ws.Range("I361").Value = "aa"
$ws.Range("J361").Value = "Agree/Accept"
$ws.Range("I374").Value = "aa"
$ws.Range("J374").Value = "Agree/Accept"
$ws.Range("I381").Value = "sd"
$ws.Range("J381").Value = "Statement-non-opinion"
$ws.Range("I384").Value = "sv"
$ws.Range("J384").Value = "Statement-opinion"
$ws.Range("I404").Value = "b"
$ws.Range("J404").Value = "Acknowledge (Backchannel)"
$ws.Range("I405").Value = "ba"
$ws.Range("J405").Value = "Appreciation"
$ws.Range("I414").Value = "b"
$ws.Range("J414").Value = "Acknowledge (Backchannel)"
$ws.Range("I416").Value = "sv"
$ws.Range("J416").Value = "Statement-opinion"
$ws.Range("I418").Value = "aa"
$ws.Range("J418").Value = "Agree/Accept"
$ws.Range("I420").Value = "aa"
$ws.Range("J420").Value = "Agree/Accept"
$ws.Range("I428").Value = "sd"
$ws.Range("J428").Value = "Statement-non-opinion"
$ws.Range("I429").Value = "sd"
$ws.Range("J429").Value = "Statement-non-opinion"
$ws.Range("I431").Value = "sd"
$ws.Range("J431").Value = "Statement-non-opinion"
$ws.Range("I445").Value = "sd"
$ws.Range("J445").Value = "Statement-non-opinion"
$ws.Range("I452").Value = "ba"
$ws.Range("J452").Value = "Appreciation"
$ws.Range("I457").Value = "b"
$ws.Range("J457").Value = "Acknowledge (Backchannel)"
$ws.Range("I458").Value = "ba"
$ws.Range("J458").Value = "Appreciation"
$ws.Range("I476").Value = "sd"
$ws.Range("J476").Value = "Statement-non-opinion"
$ws.Range("I477").Value = "sd"
$ws.Range("J477").Value = "Statement-non-opinion"
$ws.Range("I482").Value = "b"
$ws.Range("J482").Value = "Acknowledge (Backchannel)"
$ws.Range("I485").Value = "sd"
$ws.Range("J485").Value = "Statement-non-opinion"
$ws.Range("I493").Value = "sd"
$ws.Range("J493").Value = "Statement-non-opinion"
$ws.Range("I515").Value = "aa"
$ws.Range("J515").Value = "Agree/Accept"
$ws.Range("I521").Value = "sv"
$ws.Range("J521").Value = "Statement-opinion"
$ws.Range("I523").Value = "aa"
$ws.Range("J523").Value = "Agree/Accept"
$ws.Range("I531").Value = "sd"
$ws.Range("J531").Value = "Statement-non-opinion"
$ws.Range("I551").Value = "sd"
$ws.Range("J551").Value = "Statement-non-opinion"
$ws.Range("I555").Value = "qy"
$ws.Range("J555").Value = "Yes-No-Question"
$ws.Range("I558").Value = "b"
$ws.Range("J558").Value = "Acknowledge (Backchannel)"
$ws.Range("I559").Value = "ba"
$ws.Range("J559").Value = "Appreciation"
$ws.Range("I562").Value = "sv"
$ws.Range("J562").Value = "Statement-opinion"
$ws.Range("I565").Value = "sv"
$ws.Range("J565").Value = "Statement-opinion"
$ws.Range("I568").Value = "sd"
$ws.Range("J568").Value = "Statement-non-opinion"
$ws.Range("I601").Value = "sd"
$ws.Range("J601").Value = "Statement-non-opinion"
$ws.Range("I619").Value = "sv"
$ws.Range("J619").Value = "Statement-opinion"
$ws.Range("I633").Value = "%"
$ws.Range("J633").Value = "Uninterpretable"
